$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set "Experimental" value (B7) to "true" (as text, not boolean)
$ws.Range("B7").Value = "'true"

# Update the "Date" value (B8) to the new timestamp
$ws.Range("B8").Value = "2025-07-21T12:46:15+00:00"
